$d = $word.ActiveDocument

# --- Step 1: merge the split run "...food chai" + "n length. " into one run ---
$d.Content.Find.Execute(
    "on food chai" + "n length. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "on food chain length. ", 2) | Out-Null

# --- Step 2: merge "A" + " recent" + " paper by " into "A recent paper by " ---
# (this also removes the old _GoBack bookmark that sat between "A" and " recent",
#  since it will be re-created in the newly inserted paragraph below)
$d.Content.Find.Execute(
    "A recent paper by ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A recent paper by ", 2) | Out-Null

# --- Step 3: insert a brand-new paragraph before the current first paragraph ---
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(1)
$newRange = $newPara.Range

# Run 1 - plain text
$newRange.InsertAfter("This paper was previously submitted to Food Webs as ")
$run1End = $newRange.End - 1

# Run 2 - "Ms. Ref. No.:  FOOWEB-D-14-00006" (styled)
# NOTE: deliberately NOT touching Font.NameBi here - in this COM shim it
# incorrectly stamps <w:rFonts w:cs="..."/> onto every run in the
# paragraph (including the unstyled run 1), not just the target range.
$newRange.InsertAfter("Ms. Ref. No.:  FOOWEB-D-14-00006")
$run2End = $newRange.End - 1
$run2Range = $d.Range($run1End, $run2End)
$run2Range.Font.Name = "Arial"
$run2Range.Font.Size = 10
$run2Range.Font.Color = 2236962

# Run 3 - closing sentence (styled, same as run 2)
$newRange.InsertAfter(". We have revised our manuscript according to the useful suggestions of the editor and two reviewers. A more detailed description of our response to the comments on the previous version of the manuscript is included. ")
$run3End = $newRange.End - 1
$run3Range = $d.Range($run2End, $run3End)
$run3Range.Font.Name = "Arial"
$run3Range.Font.Size = 10
$run3Range.Font.Color = 2236962

# bookmark "_GoBack" at the end of the new paragraph.
# Placing a bookmark exactly at "end of paragraph, right before the
# paragraph mark" is ambiguous when a following paragraph exists (it
# snaps to the start of the next paragraph instead). Work around this
# by appending a throw-away sentinel character, anchoring the bookmark
# right before it (now a safe, unambiguous mid-run position), and then
# deleting the sentinel.
$newRange.InsertAfter("X")
$sentinelEnd = $newRange.End - 1
$bmPos = $sentinelEnd - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

$newPara = $d.Paragraphs(1)
$sentinelRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$sentinelRange.Delete() | Out-Null
